$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the currently-empty row 2 cells (A2:D2) so the row becomes truly empty,
# but keep the row itself present (touching OutlineLevel forces the row element to persist).
$ws.Range("A2:D2").ClearContents()
$ws.Rows(2).OutlineLevel = 0

# Populate new row 3 with appointment data (kept as text, matching inlineStr source)
$ws.Range("A3:D3").NumberFormat = "@"
$ws.Range("A3").Value = "lunes"
$ws.Range("B3").Value = "10"
$ws.Range("C3").Value = "maría"
$ws.Range("D3").Value = "lópez"
